$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: capacitor C2,C3,C4,C5,C6 quantity increased from 3 to 5
$ws.Range("D3").Value = 5
$ws.Range("C3").Value = "C2, C3, C4, C5, C6"

# Row 16: remove the ON Schottky Diode line entirely (was D2 part)
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Row -eq 16) {
        $hl.Delete()
    }
}
$ws.Range("A16:H16").ClearContents()

# Column C widened to fit the longer designator list
$ws.Columns.Item(3).AutoFit()

# Selection moved (per commit: cursor left near bottom of sheet after edits)
$ws.Range("B24").Select() | Out-Null
